# Generate Report for Handback
#
# Updates the localization-status workbook to reflect that the zh-cn and
# de-de handback packages have come back "in sync with en-US":
#   - Status cells move from "Ready for handoff" -> "Handed back: in sync
#     with en-US" (shared across the Overview summary + both language
#     sheets because they point at the same shared string).
#   - Each language sheet's "Latest Target File" (I2) / "Latest Handback
#     File" (J2) columns get populated with the handed-back filenames, and
#     I2 becomes a hyperlink to the source markdown file (matching the
#     style already used for column A).
#   - "Latest Handback DateTime" (K2) is stamped with the handback time.
#   - The now much-longer filename values mean columns C/I/J (and the
#     Overview sheet's E/F) need to be widened so the text is not clipped.

$wb = $excel.ActiveWorkbook

$mdFile = "899ee086-8242-4535-95c3-0bab5ea32bdf.md"
$mdUrl  = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/d72acab8a6643c44e28d3405219577ae693d83aa/e2e/$mdFile"
$newStatus = "Handed back: in sync with en-US"

# --- Overview sheet: status summary columns (E = zh-cn, F = de-de) -----
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E2").Value = $newStatus
$overview.Range("F2").Value = $newStatus
$overview.Columns.Item(5).ColumnWidth = 29.9777047293527
$overview.Columns.Item(6).ColumnWidth = 29.9777047293527

# --- Per-language sheets -------------------------------------------------
$languages = @(
    @{ Sheet = "zh-cn"; Handback = "899ee086-8242-4535-95c3-0bab5ea32bdf.ee505f28e96fef06e3947d82f0999e12df38f31f.zh-cn.xlf"; DateTime = "2016-08-15 22:55:13" },
    @{ Sheet = "de-de"; Handback = "899ee086-8242-4535-95c3-0bab5ea32bdf.ee505f28e96fef06e3947d82f0999e12df38f31f.de-de.xlf"; DateTime = "2016-08-15 22:55:20" }
)

foreach ($lang in $languages) {
    $ws = $wb.Worksheets.Item($lang.Sheet)

    # Status (Column C)
    $ws.Range("C2").Value = $newStatus

    # Latest Target File (Column I) - new hyperlink to the source md file
    $ws.Hyperlinks.Add($ws.Range("I2"), $mdUrl, "", "", $mdFile)
    $ws.Range("I2").Font.Underline = 2
    $ws.Range("I2").Font.Color = 15570276

    # Latest Handback File (Column J)
    $ws.Range("J2").Value = $lang.Handback

    # Latest Handback DateTime (Column K)
    $ws.Range("K2").Value = $lang.DateTime

    # Widen columns C, I, J to fit the longer text values
    $ws.Columns.Item(3).ColumnWidth = 29.9777047293527
    $ws.Columns.Item(9).ColumnWidth = 40
    $ws.Columns.Item(10).ColumnWidth = 40
}

Write-Host "Handback report generated"
